# facture_txt update:
# The final paragraph ("Total général : 0 €.") is replaced by three new
# "line item" paragraphs (Bonbons / Biscuits / Laits, each with a manual
# line break before the quantity/price/total detail) followed by an
# updated grand-total paragraph ("Total général : 1000 €.").

$d = $word.ActiveDocument

# Locate the paragraph that holds the old "Total général : 0 €." line by
# content (rather than assuming a fixed index) and replace it in place via
# an OOXML fragment so we get the exact <w:t>/<w:br/>/<w:t> run structure.
$searchRange = $d.Content
$found = $searchRange.Find.Execute("Total général : 0", $true, $false, $false,
                                    $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $targetParagraph = $searchRange.Paragraphs.First
} else {
    $targetParagraph = $d.Paragraphs.Last
}
$targetRange = $targetParagraph.Range

$newParagraphsXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Bonbons, </w:t><w:br/><w:t xml:space="preserve"> - Quantité : 2, Prix unitaire : 150€, Total : 300€</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Biscuits, </w:t><w:br/><w:t xml:space="preserve"> - Quantité : 2, Prix unitaire : 200€, Total : 400€</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Laits, </w:t><w:br/><w:t xml:space="preserve"> - Quantité : 3, Prix unitaire : 100€, Total : 300€</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Total général : 1000 €.</w:t></w:r></w:p>'

$targetRange.InsertXML($newParagraphsXml)
